$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.069.44"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "1.836.77"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.50"
$ws.Range("E5").Value = "  -2.91%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3873"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07865"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9624"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.92"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").Value = "1.848.88"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.700"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.929"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06857"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.60"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009953"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.70"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "28.074.24"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.333"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.098"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "2.046.55"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.43"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.17"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.684"
$ws.Range("E28").Value = "  -7.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.964"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.29"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9391"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09237"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.280"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.324"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.307"
$ws.Range("E35").Value = "  -4.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05863"
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02131"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.139"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.796"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5598"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.916"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1765"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07246"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.71"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5273"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.125"
$ws.Range("E46").Value = "  -10.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.120"
$ws.Range("E47").Value = "  -10.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.832"
$ws.Range("E48").Value = "  -4.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.63"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.030"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  +0.05%  "
